# This script updates the "LR-pairs_lrc2p/Efnb3-Epha4" results sheet with
# freshly recomputed TPM-based NATMI statistics, and adds a new sending
# cluster ("Resolving-Mac") that was produced by re-running the analysis
# with the updated cell-type annotations (a prior "Resolving-Mac" cluster
# was renamed to "Inflammatory-Mac").
#
# Columns A:T are, in order:
#   A Sending cluster            K Receptor-expressing cells
#   B Ligand symbol              L Receptor detection rate
#   C Receptor symbol            M Receptor average expression value
#   D Target cluster             N Receptor total expression value
#   E Ligand-expressing cells    O Receptor derived specificity (avg)
#   F Ligand detection rate      P Receptor derived specificity (total)
#   G Ligand average expr value  Q Edge average expression weight
#   H Ligand total expr value    R Edge total expression weight
#   I Ligand derived spec (avg)  S Edge average expression derived spec
#   J Ligand derived spec (tot)  T Edge total expression derived spec

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efnb3", "Epha4", "ECs", 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 3, 1, 8.540560666666666, 25.621682, 0.4159358086620884, 0.4159358086620884, 1.309410292877778, 11.7846926359, 0.04582283035059939, 0.04582283035059938),
    @("ECs", "Efnb3", "Epha4", "FAPs", 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 3, 1, 11.57455166666667, 34.723655, 0.563694901924408, 0.563694901924408, 1.774571679694445, 15.97114511725, 0.06210115917517602, 0.06210115917517602),
    @("ECs", "Efnb3", "Epha4", "Inflammatory-Mac", 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 1, 0.3333333333333333, 0.026642, 0.079926, 0.001297498167494471, 0.001297498167494471, 0.004084662633333334, 0.0367619637, 0.0001429428223565497, 0.0001429428223565497),
    @("ECs", "Efnb3", "Epha4", "MuSCs", 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 3, 1, 0.391608, 1.174824, 0.01907179124600912, 0.01907179124600912, 0.06004003320000001, 0.5403602988, 0.002101101748269789, 0.002101101748269789),
    @("FAPs", "Efnb3", "Epha4", "ECs", 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 3, 1, 8.540560666666666, 25.621682, 0.4159358086620884, 0.4159358086620884, 8.164741835090668, 73.48267651581601, 0.2857252474650611, 0.2857252474650611),
    @("FAPs", "Efnb3", "Epha4", "FAPs", 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 3, 1, 11.57455166666667, 34.723655, 0.563694901924408, 0.563694901924408, 11.06522509512667, 99.58702585614002, 0.3872276971420693, 0.3872276971420692),
    @("FAPs", "Efnb3", "Epha4", "Inflammatory-Mac", 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 1, 0.3333333333333333, 0.026642, 0.079926, 0.001297498167494471, 0.001297498167494471, 0.025469645432, 0.229226808888, 0.0008913105755075907, 0.0008913105755075906),
    @("FAPs", "Efnb3", "Epha4", "MuSCs", 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 3, 1, 0.391608, 1.174824, 0.01907179124600912, 0.01907179124600912, 0.3743756815680001, 3.369381134112001, 0.01310128188024084, 0.01310128188024084),
    @("MuSCs", "Efnb3", "Epha4", "ECs", 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 3, 1, 8.540560666666666, 25.621682, 0.4159358086620884, 0.4159358086620884, 1.080494798475556, 9.72445318628, 0.03781192962553863, 0.03781192962553863),
    @("MuSCs", "Efnb3", "Epha4", "FAPs", 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 3, 1, 11.57455166666667, 34.723655, 0.563694901924408, 0.563694901924408, 1.464335113188889, 13.1790160187, 0.05124442646667313, 0.05124442646667313),
    @("MuSCs", "Efnb3", "Epha4", "Inflammatory-Mac", 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 1, 0.3333333333333333, 0.026642, 0.079926, 0.001297498167494471, 0.001297498167494471, 0.003370568226666667, 0.03033511404, 0.0001179530792416673, 0.0001179530792416673),
    @("MuSCs", "Efnb3", "Epha4", "MuSCs", 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 3, 1, 0.391608, 1.174824, 0.01907179124600912, 0.01907179124600912, 0.04954363344000001, 0.44589270096, 0.001733780101181249, 0.001733780101181249),
    @("Resolving-Mac", "Efnb3", "Epha4", "ECs", 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 3, 1, 8.540560666666666, 25.621682, 0.4159358086620884, 0.4159358086620884, 1.330926812050667, 11.978341308456, 0.04657580122088927, 0.04657580122088927),
    @("Resolving-Mac", "Efnb3", "Epha4", "FAPs", 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 3, 1, 11.57455166666667, 34.723655, 0.563694901924408, 0.563694901924408, 1.803731833526667, 16.23358650174, 0.0631216191404896, 0.06312161914048958),
    @("Resolving-Mac", "Efnb3", "Epha4", "Inflammatory-Mac", 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 1, 0.3333333333333333, 0.026642, 0.079926, 0.001297498167494471, 0.001297498167494471, 0.004151782712, 0.037366044408, 0.0001452916903886636, 0.0001452916903886636),
    @("Resolving-Mac", "Efnb3", "Epha4", "MuSCs", 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 3, 1, 0.391608, 1.174824, 0.01907179124600912, 0.01907179124600912, 0.061026624288, 0.549239618592, 0.002135627516317235, 0.002135627516317235)
)

$startRow = 2
$r = $startRow
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
